$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2458.8
$ws.Range("I40").Value = 2150
$ws.Range("J40").Value = 2506.3076
$ws.Range("K40").Value = 2150
$ws.Range("L40").Value = 2506.3076
$ws.Range("M40").Value = -1975
$ws.Range("N40").Value = -2856.3076

$ws.Range("H70").Value = 1473.7407
$ws.Range("I70").Value = 1523.5883
$ws.Range("J70").Value = 1389
$ws.Range("K70").Value = 4570.7649
$ws.Range("L70").Value = 4167
$ws.Range("M70").Value = -4300.7649
$ws.Range("N70").Value = -4707

$ws.Range("H73").Value = 1473.7407
$ws.Range("I73").Value = 1523.5883
$ws.Range("J73").Value = 1389
$ws.Range("K73").Value = 4570.7649
$ws.Range("L73").Value = 4167
$ws.Range("M73").Value = -3634.7649
$ws.Range("N73").Value = -6039

$ws.Range("H108").Value = 39759.5
$ws.Range("J108").Value = 39759.5
$ws.Range("L108").Value = 39759.5
$ws.Range("N108").Value = -47439.5

$ws.Range("H129").Value = 2580.2415
$ws.Range("I129").Value = 13005.5
$ws.Range("J129").Value = 912.2
$ws.Range("K129").Value = 39016.5
$ws.Range("L129").Value = 2736.6
$ws.Range("M129").Value = -34016.5
$ws.Range("N129").Value = -12736.6

$ws.Range("H138").Value = 2977.7412
$ws.Range("I138").Value = 1822.68
$ws.Range("J138").Value = 3459.0166
$ws.Range("K138").Value = 5468.04
$ws.Range("L138").Value = 10377.0498
$ws.Range("M138").Value = -328.04
$ws.Range("N138").Value = -20657.0498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 5700
$ws.Range("J44").Value = 5700
$ws.Range("L44").Value = 5700
$ws.Range("N44").Value = -6676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 169.1
$ws.Range("I22").Value = 180.125
$ws.Range("J22").Value = 125
$ws.Range("K22").Value = 180.125
$ws.Range("L22").Value = 125
$ws.Range("M22").Value = -7.125
$ws.Range("N22").Value = -471

$ws.Range("H86").Value = 222121.2
$ws.Range("I86").Value = 277251.5
$ws.Range("K86").Value = 277251.5
$ws.Range("M86").Value = -276128.5

$ws.Range("H89").Value = 222121.2
$ws.Range("I89").Value = 277251.5
$ws.Range("K89").Value = 1386257.5
$ws.Range("M89").Value = -1380641.5

$ws.Range("H107").Value = 100046790
$ws.Range("I107").Value = 125058104
$ws.Range("J107").Value = 1550
$ws.Range("K107").Value = 125058104
$ws.Range("L107").Value = 1550
$ws.Range("M107").Value = -125056184
$ws.Range("N107").Value = -5390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 23700
$ws.Range("J29").Value = 23700
$ws.Range("L29").Value = 23700
$ws.Range("N29").Value = -24286

$ws.Range("H31").Value = 58018.73
$ws.Range("I31").Value = 1378.091
$ws.Range("J31").Value = 99555.2
$ws.Range("K31").Value = 1378.091
$ws.Range("L31").Value = 99555.2
$ws.Range("M31").Value = -1083.091
$ws.Range("N31").Value = -100145.2

$ws.Range("H34").Value = 58018.73
$ws.Range("I34").Value = 1378.091
$ws.Range("J34").Value = 99555.2
$ws.Range("K34").Value = 1378.091
$ws.Range("L34").Value = 99555.2
$ws.Range("M34").Value = -1176.091
$ws.Range("N34").Value = -99959.2

$ws.Range("H50").Value = 9500
$ws.Range("J50").Value = 9500
$ws.Range("L50").Value = 9500
$ws.Range("N50").Value = -10750

$ws.Range("H51").Value = 7897.375
$ws.Range("J51").Value = 7897.375
$ws.Range("L51").Value = 7897.375
$ws.Range("N51").Value = -9369.375

$ws.Range("H58").Value = 1749.1875
$ws.Range("I58").Value = 1526.25
$ws.Range("J58").Value = 2418
$ws.Range("K58").Value = 1526.25
$ws.Range("L58").Value = 2418
$ws.Range("M58").Value = -1323.25
$ws.Range("N58").Value = -2824

$ws.Range("H60").Value = 14612.6
$ws.Range("J60").Value = 14612.6
$ws.Range("L60").Value = 14612.6
$ws.Range("N60").Value = -15634.6

$ws.Range("H61").Value = 7897.375
$ws.Range("J61").Value = 7897.375
$ws.Range("L61").Value = 7897.375
$ws.Range("N61").Value = -8593.375

$ws.Range("H74").Value = 38448.57
$ws.Range("J74").Value = 38448.57
$ws.Range("L74").Value = 38448.57
$ws.Range("N74").Value = -40196.57

$ws.Range("H77").Value = 38448.57
$ws.Range("J77").Value = 38448.57
$ws.Range("L77").Value = 115345.71
$ws.Range("N77").Value = -124081.71

$ws.Range("H122").Value = 981.5714
$ws.Range("I122").Value = 796.6667
$ws.Range("J122").Value = 1120.25
$ws.Range("K122").Value = 2390.0001
$ws.Range("L122").Value = 3360.75
$ws.Range("M122").Value = 59.9998999999998
$ws.Range("N122").Value = -8260.75

$ws.Range("H132").Value = 2070.8147
$ws.Range("I132").Value = 1975.2439
$ws.Range("J132").Value = 2372.2307
$ws.Range("K132").Value = 5925.7317
$ws.Range("L132").Value = 7116.6921
$ws.Range("M132").Value = -3395.7317
$ws.Range("N132").Value = -12176.6921

$ws.Range("H136").Value = 1749.1875
$ws.Range("I136").Value = 1526.25
$ws.Range("J136").Value = 2418
$ws.Range("K136").Value = 4578.75
$ws.Range("L136").Value = 7254
$ws.Range("M136").Value = -2028.75
$ws.Range("N136").Value = -12354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1176.5491
$ws.Range("I5").Value = 792.0769
$ws.Range("J5").Value = 1576.4
$ws.Range("K5").Value = 2376.2307
$ws.Range("L5").Value = 4729.200000000001
$ws.Range("M5").Value = -2264.2307
$ws.Range("N5").Value = -4953.200000000001

$ws.Range("H37").Value = 610829.25
$ws.Range("J37").Value = 610829.25
$ws.Range("L37").Value = 1832487.75
$ws.Range("N37").Value = -1832711.75

$ws.Range("H55").Value = 11056.363
$ws.Range("J55").Value = 3308.3333
$ws.Range("L55").Value = 9924.999899999999
$ws.Range("N55").Value = -10278.9999

$ws.Range("H131").Value = 900.0700000000001
$ws.Range("J131").Value = 962.9167
$ws.Range("L131").Value = 2888.7501
$ws.Range("N131").Value = -12968.7501

$ws.Range("H135").Value = 1176.5491
$ws.Range("I135").Value = 792.0769
$ws.Range("J135").Value = 1576.4
$ws.Range("K135").Value = 7128.6921
$ws.Range("L135").Value = 14187.6
$ws.Range("M135").Value = -4593.6921
$ws.Range("N135").Value = -19257.6

$ws.Range("H137").Value = 2734319.8
$ws.Range("I137").Value = 68564.60000000001
$ws.Range("J137").Value = 4551880
$ws.Range("K137").Value = 205693.8
$ws.Range("L137").Value = 13655640
$ws.Range("M137").Value = -200593.8
$ws.Range("N137").Value = -13665840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2767
$ws.Range("I132").Value = 1971.2858
$ws.Range("J132").Value = 4052.3845
$ws.Range("K132").Value = 5913.857400000001
$ws.Range("L132").Value = 12157.1535
$ws.Range("M132").Value = -3383.857400000001
$ws.Range("N132").Value = -17217.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2555.25
$ws.Range("I7").Value = 1777.1666
$ws.Range("J7").Value = 3333.3333
$ws.Range("K7").Value = 1777.1666
$ws.Range("L7").Value = 3333.3333
$ws.Range("M7").Value = -1665.1666
$ws.Range("N7").Value = -3557.3333

$ws.Range("H16").Value = 7161077
$ws.Range("I16").Value = 9695305
$ws.Range("J16").Value = 1670249.9
$ws.Range("K16").Value = 9695305
$ws.Range("L16").Value = 1670249.9
$ws.Range("M16").Value = -9695135
$ws.Range("N16").Value = -1670589.9

$ws.Range("H126").Value = 2555.25
$ws.Range("I126").Value = 1777.1666
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 5331.4998
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -2861.4998
$ws.Range("N126").Value = -14939.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2159.25
$ws.Range("I126").Value = 2032
$ws.Range("K126").Value = 6096
$ws.Range("M126").Value = -3626
